# Add two new history rows (17 and 18) to the "historique" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: Sortie
$ws.Cells.Item(17, 1).Value = "2025-05-23 07:59:24"
$ws.Cells.Item(17, 2).Value = "Perceuse sans fil"
$ws.Cells.Item(17, 3).Value = "Sortie"
$ws.Cells.Item(17, 4).Value = 6
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 3

# Row 18: Entrée
$ws.Cells.Item(18, 1).Value = "2025-05-23 08:09:59"
$ws.Cells.Item(18, 2).Value = "Perceuse sans fil"
$ws.Cells.Item(18, 3).Value = "Entrée"
$ws.Cells.Item(18, 4).Value = 8
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 11
